# Generate Report for Handoff
# Adds two new source files (613ed53d-... and eb53d95b-...) to the
# localization-status workbook: one new row on "Overview", and one new row
# each on the "zh-cn" / "de-de" per-language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet (File Name / Path And Name / Extension / Publish URL /
#                 zh-cn / de-de / Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$overviewRows = @(
    @("613ed53d-8f16-4b81-8864-9ec5421e3b62.md", "e2e\613ed53d-8f16-4b81-8864-9ec5421e3b62.md", "Ready for handoff", "2016-09-01 00:45:16"),
    @("eb53d95b-d256-41f4-88b9-8b0024a2ca90.md", "e2e\eb53d95b-d256-41f4-88b9-8b0024a2ca90.md", "Ready for handoff", "2016-09-01 00:45:16")
)

foreach ($data in $overviewRows) {
    $fileName = $data[0]
    $pathName = $data[1]
    $status   = $data[2]
    $date     = $data[3]

    $newRow = $loOverview.ListRows.Add()
    $r = $newRow.Range

    $r.Cells.Item(1, 1).Value2 = $fileName
    $r.Cells.Item(1, 2).Value2 = $pathName
    $r.Cells.Item(1, 3).Value2 = ".md"
    $r.Cells.Item(1, 4).Value2 = ""
    $r.Cells.Item(1, 5).Value2 = $status
    $r.Cells.Item(1, 6).Value2 = $status
    $r.Cells.Item(1, 7).Value2 = $date

    $wsOverview.Hyperlinks.Add($r.Cells.Item(1, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/$([System.IO.Path]::GetFileName($fileName))", "", "", $pathName) | Out-Null
}

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn / de-de) — 16 columns:
# Source File Name, File Extension, Status, Source Path, Priority,
# Content Duplicate, Latest Handoff File, Latest Handoff Datetime,
# Latest Target File, Latest Handback File, Latest Handback DateTime,
# Reference Tokens, To be localized, Dependency From, Has metadata,
# Error Detail
# ---------------------------------------------------------------------
$languages = @(
    @{
        Sheet = "zh-cn"
        RepoSuffix = "ol-test0-zhcn"
        Rows = @(
            @{
                SourceFileName = "613ed53d-8f16-4b81-8864-9ec5421e3b62.md"
                HandoffFile    = "613ed53d-8f16-4b81-8864-9ec5421e3b62.23fa9e4e2a683e045a41530aef4f877b1a9bf941.zh-cn.xlf"
                HandoffDate    = "2016-09-01 00:45:08"
            },
            @{
                SourceFileName = "eb53d95b-d256-41f4-88b9-8b0024a2ca90.md"
                HandoffFile    = "eb53d95b-d256-41f4-88b9-8b0024a2ca90.62276ca71ab1c7cc0f3db7234859d6554aaad307.zh-cn.xlf"
                HandoffDate    = "2016-09-01 00:45:08"
            }
        )
    },
    @{
        Sheet = "de-de"
        RepoSuffix = "ol-test0-dede"
        Rows = @(
            @{
                SourceFileName = "613ed53d-8f16-4b81-8864-9ec5421e3b62.md"
                HandoffFile    = "613ed53d-8f16-4b81-8864-9ec5421e3b62.23fa9e4e2a683e045a41530aef4f877b1a9bf941.de-de.xlf"
                HandoffDate    = "2016-09-01 00:45:16"
            },
            @{
                SourceFileName = "eb53d95b-d256-41f4-88b9-8b0024a2ca90.md"
                HandoffFile    = "eb53d95b-d256-41f4-88b9-8b0024a2ca90.62276ca71ab1c7cc0f3db7234859d6554aaad307.de-de.xlf"
                HandoffDate    = "2016-09-01 00:45:16"
            }
        )
    }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)
    $lo = $ws.ListObjects.Item(1)

    foreach ($row in $lang.Rows) {
        $newRow = $lo.ListRows.Add()
        $r = $newRow.Range

        $r.Cells.Item(1, 1).Value2  = $row.SourceFileName   # Source File Name
        $r.Cells.Item(1, 2).Value2  = ".md"                 # File Extension
        $r.Cells.Item(1, 3).Value2  = "Ready for handoff"    # Status
        $r.Cells.Item(1, 4).Value2  = "e2e"                 # Source Path
        $r.Cells.Item(1, 5).Value2  = "ht"                  # Priority
        $r.Cells.Item(1, 6).Value2  = "'False"              # Content Duplicate (force text, not boolean)
        $r.Cells.Item(1, 7).Value2  = $row.HandoffFile       # Latest Handoff File
        $r.Cells.Item(1, 8).Value2  = $row.HandoffDate       # Latest Handoff Datetime
        $r.Cells.Item(1, 9).Value2  = ""                     # Latest Target File
        $r.Cells.Item(1, 10).Value2 = ""                     # Latest Handback File
        $r.Cells.Item(1, 11).Value2 = "0001-01-01 00:00:00"  # Latest Handback DateTime
        $r.Cells.Item(1, 12).Value2 = ""                     # Reference Tokens
        $r.Cells.Item(1, 13).Value2 = "'True"               # To be localized (force text, not boolean)
        $r.Cells.Item(1, 14).Value2 = ""                     # Dependency From
        $r.Cells.Item(1, 15).Value2 = "'False"              # Has metadata (force text, not boolean)
        $r.Cells.Item(1, 16).Value2 = ""                     # Error Detail

        $target = "https://github.com/OpenLocalizationTestOrg/$($lang.RepoSuffix)/blob/HEAD/e2e/$($row.SourceFileName)"
        $ws.Hyperlinks.Add($r.Cells.Item(1, 1), $target, "", "", $row.SourceFileName) | Out-Null
    }
}
